# Week 10 Lab - Social Support
#
# The final table's "Support 3" answer cell is currently empty (its
# paragraph holds nothing but a manual line break, <w:br/>). The edit
# fills that answer in with three runs of text (all sharing the same
# "Helvetica Light" / sz 22 formatting used throughout the table).
#
# We locate the target paragraph via the table/cell indices (last table,
# row 3 "Support 3", column 2) rather than Find, since the cell's only
# content is a non-text line break and there is no text anchor to search
# for. We then use Range.InsertXML on a range that spans just that break
# character so the paragraph's existing identity (w14:paraId/rsid/pPr) is
# preserved while its run content is replaced by the new runs.

$d = $word.ActiveDocument

# The document has four identically-shaped "Support 1/2/3" tables; the
# one we need is the last one, and within it the last row ("Support 3").
$table = $d.Tables.Item($d.Tables.Count)
$lastRow = $table.Rows.Count
$labelCell = $table.Cell($lastRow, 1)
$cell = $table.Cell($lastRow, 2)
$cellRange = $cell.Range

$labelText = $labelCell.Range.Text.TrimEnd([char]13, [char]7).Trim()
if ($labelText -ne "Support 3") {
    throw "Expected last table's last row to be labelled 'Support 3', found: " + $labelText
}

# Target just the line-break character inside the paragraph (not the
# paragraph end mark / cell mark that follow it).
$target = $d.Range($cellRange.Start, $cellRange.Start + 1)

$openXmlTemplate = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="67A212EB" w14:textId="3CDA5993" w:rsidR="005C285B" w:rsidRDefault="003B3939" w:rsidP="00B55B41">
            <w:pPr>
              <w:spacing w:line="276" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Arial"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Arial"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
              </w:rPr>
              <w:t xml:space="preserve">__RUN1__</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Arial"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
              </w:rPr>
              <w:t xml:space="preserve">__RUN2__</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Arial"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
              </w:rPr>
              <w:t xml:space="preserve">__RUN3__</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$run1 = "There is not much to do to make the situation better as I go down to Toronto every day that I don" + [char]0x2019 + "t have to go into the office for work. I spend the maximum amount of time I can each week down with her. A step I could take is conversing about meal prep ideas for multiple weeks in a row, so "
$run2 = "the process "
$run3 = "is as easy as possible. "

$openXml = $openXmlTemplate.Replace("__RUN1__", $run1).Replace("__RUN2__", $run2).Replace("__RUN3__", $run3)

$target.InsertXML($openXml)

$updatedCell = $table.Cell($lastRow, 2)
Write-Host "Support 3 answer now reads:" $updatedCell.Range.Text
